$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 83439
$ws.Range("B2").Value = "Ana Sophia Martins"
$ws.Range("D2").Value = "Doença"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45102
$ws.Range("G2").Value = 11872.35

# Row 3
$ws.Range("A3").Value = 44412
$ws.Range("B3").Value = "Brenda Correia"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45105
$ws.Range("G3").Value = 5794.94

# Row 4
$ws.Range("A4").Value = 37584
$ws.Range("B4").Value = "Lavínia Almeida"
$ws.Range("C4").Value = "Atendimento ao Cliente"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 45082
$ws.Range("G4").Value = 6774.65

# Row 5
$ws.Range("A5").Value = 2821
$ws.Range("B5").Value = "Alexia Cardoso"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45090
$ws.Range("G5").Value = 3943.37

# Row 6
$ws.Range("A6").Value = 6810
$ws.Range("B6").Value = "Sra. Ana Beatriz Campos"
$ws.Range("C6").Value = "Recursos Humanos"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 45087
$ws.Range("G6").Value = 12079.1

# Row 7
$ws.Range("A7").Value = 98571
$ws.Range("B7").Value = "Davi Lucas Mendes"
$ws.Range("C7").Value = "TI"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45104
$ws.Range("G7").Value = 12344.61

# Row 8
$ws.Range("A8").Value = 96249
$ws.Range("B8").Value = "Sr. Cauê Araújo"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45084
$ws.Range("G8").Value = 12015.67

# Row 9
$ws.Range("A9").Value = 34266
$ws.Range("B9").Value = "Noah Correia"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45089
$ws.Range("G9").Value = 4883.38

# Row 10
$ws.Range("A10").Value = 63204
$ws.Range("B10").Value = "Bernardo Lopes"
$ws.Range("C10").Value = "Recursos Humanos"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45104
$ws.Range("G10").Value = 10240.39

# Row 11
$ws.Range("A11").Value = 24020
$ws.Range("B11").Value = "Alexia Almeida"
$ws.Range("C11").Value = "P&D"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45095
$ws.Range("G11").Value = 9781.799999999999
